$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# "prix" (col C) and "vitesse" (col D) used to be placeholder values (0 / 1).
# Fill them in with real per-km price and average speed figures, row-by-row,
# for each of the 19 transport modes (rows 2-20). "mode_transport" (col A)
# and "emissions" (col B) are untouched.
# ---------------------------------------------------------------------------
$prixVitesse = @(
    @(0.15, 600),   # 2  Avion (court-courrier)
    @(0.13, 80),    # 3  Moto
    @(0.1,  800),   # 4  Avion (moyen-courrier)
    @(0.08, 900),   # 5  Avion (long-courrier)
    @(0.33, 100),   # 6  Voiture thermique
    @(0.1,  50),    # 7  Bus thermique
    @(0.13, 45),    # 8  Scooter et moto légère
    @(0.27, 100),   # 9  Voiture électrique
    @(0.15, 90),    # 10 TER
    @(0.05, 80),    # 11 Autocar
    @(0.05, 25),    # 12 Vélo à assistance électrique
    @(0.02, 20),    # 13 Trottinette électrique
    @(0.12, 50),    # 14 RER ou Transilien
    @(0.1,  100),   # 15 Train Intercités
    @(0.12, 30),    # 16 Métro
    @(0.1,  20),    # 17 Tramway
    @(0.15, 300),   # 18 TGV
    @(0,    15),    # 19 Vélo
    @(0,    5)      # 20 Marche
)

for ($i = 0; $i -lt $prixVitesse.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $prixVitesse[$i][0]
    $ws.Cells.Item($row, 4).Value = $prixVitesse[$i][1]
}

# ---------------------------------------------------------------------------
# Formatting
# ---------------------------------------------------------------------------
# Header row now wraps text (vertical alignment left at the default).
$ws.Range("A1:D1").WrapText = $true

# Column D (vitesse) data cells pick up the same wrap + vertical-center
# formatting already used by columns B (emissions) and C (prix).
$ws.Range("D2:D20").WrapText = $true
$ws.Range("D2:D20").VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# Column widths
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 23.5
$ws.Columns.Item(2).ColumnWidth = 28
$ws.Columns.Item(3).ColumnWidth = 37
$ws.Columns.Item(4).ColumnWidth = 30.5

# ---------------------------------------------------------------------------
# View state: zoom level + selected cell
# ---------------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 63
$ws.Range("K11").Select() | Out-Null
